$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2; existing rows 2-6 shift down to 3-7.
$ws.Rows.Item(2).Insert(-4121)  # -4121 = xlShiftDown

# Row-insert copies the header row formatting onto the new row; the source
# data rows carry no explicit style, so strip it back off.
$ws.Range("A2:AO2").ClearFormats()

# Column B ("Date") is literal text "2025-11-12" on every row and must stay
# text (not auto-converted to a date serial). Copy it from the row below,
# which already carries that exact literal-text cell after the shift.
$ws.Range("B3").Copy($ws.Range("B2")) | Out-Null

# --- Row 2 ---
$ws.Range("A2").Value = "Friendly Matches"
$ws.Range("C2").Value = "08:00:00"
$ws.Range("D2").Value = "Hamburger SV"
$ws.Range("E2").Value = "FC Groningen"
$row2 = New-Object 'object[,]' 1,36
$row2[0,0] = 2.24
$row2[0,1] = 2.4
$row2[0,2] = 2.94
$row2[0,3] = 3.25
$row2[0,4] = 3.8
$row2[0,5] = 4.5
$row2[0,6] = 1.29
$row2[0,7] = 1.03
$row2[0,8] = 5.5
$row2[0,9] = 1.18
$row2[0,10] = 2.54
$row2[0,11] = 1.58
$row2[0,12] = 1.62
$row2[0,13] = 2.38
$row2[0,14] = 1.52
$row2[0,15] = 2.58
$row2[0,16] = 1.45
$row2[0,17] = 1.72
$row2[0,18] = 980
$row2[0,19] = 980
$row2[0,20] = 980
$row2[0,21] = 220
$row2[0,22] = 980
$row2[0,23] = 42
$row2[0,24] = 980
$row2[0,25] = 980
$row2[0,26] = 980
$row2[0,27] = 980
$row2[0,28] = 980
$row2[0,29] = 980
$row2[0,30] = 980
$row2[0,31] = 980
$row2[0,32] = 980
$row2[0,33] = 580
$row2[0,34] = 600
$row2[0,35] = 980
$ws.Range("F2:AO2").Value = $row2

# --- Row 3 ---
$ws.Range("A3").Value = "Friendly Matches"
$ws.Range("C3").Value = "13:00:00"
$ws.Range("D3").Value = "SV Lafnitz"
$ws.Range("E3").Value = "KSV 1919"
$row3 = New-Object 'object[,]' 1,36
$row3[0,0] = 3.6
$row3[0,1] = 4.5
$row3[0,2] = 1.7
$row3[0,3] = 1.84
$row3[0,4] = 4.3
$row3[0,5] = 6.6
$row3[0,6] = 1.16
$row3[0,7] = 1.02
$row3[0,8] = 9.4
$row3[0,9] = 1.08
$row3[0,10] = 4.4
$row3[0,11] = 1.24
$row3[0,12] = 2.4
$row3[0,13] = 1.61
$row3[0,14] = 1.35
$row3[0,15] = 3.3
$row3[0,16] = 2.16
$row3[0,17] = 1.29
$row3[0,18] = 1000
$row3[0,19] = 1000
$row3[0,20] = 1000
$row3[0,21] = 1000
$row3[0,22] = 1000
$row3[0,23] = 1000
$row3[0,24] = 1000
$row3[0,25] = 1000
$row3[0,26] = 1000
$row3[0,27] = 1000
$row3[0,28] = 1000
$row3[0,29] = 1000
$row3[0,30] = 1000
$row3[0,31] = 1000
$row3[0,32] = 1000
$row3[0,33] = 1000
$row3[0,34] = 1000
$row3[0,35] = 15
$ws.Range("F3:AO3").Value = $row3

# --- Row 4 ---
$ws.Range("A4").Value = "Bosnian Premier League"
$ws.Range("C4").Value = "14:00:00"
$ws.Range("D4").Value = "Borac Banja Luka"
$ws.Range("E4").Value = "Zrinjski"
$row4 = New-Object 'object[,]' 1,36
$row4[0,0] = 2.48
$row4[0,1] = 2.88
$row4[0,2] = 3.1
$row4[0,3] = 3.5
$row4[0,4] = 3
$row4[0,5] = 3.5
$row4[0,6] = 1.36
$row4[0,7] = 1.1
$row4[0,8] = 2.98
$row4[0,9] = 1.3
$row4[0,10] = 1.25
$row4[0,11] = 1.35
$row4[0,12] = 1.22
$row4[0,13] = 2
$row4[0,14] = 1.05
$row4[0,15] = 1.04
$row4[0,16] = 1.4
$row4[0,17] = 1.54
$row4[0,18] = 1000
$row4[0,19] = 1000
$row4[0,20] = 42
$row4[0,21] = 65
$row4[0,22] = 1000
$row4[0,23] = 42
$row4[0,24] = 15.5
$row4[0,25] = 48
$row4[0,26] = 24
$row4[0,27] = 12
$row4[0,28] = 21
$row4[0,29] = 65
$row4[0,30] = 42
$row4[0,31] = 34
$row4[0,32] = 55
$row4[0,33] = 580
$row4[0,34] = 29
$row4[0,35] = 600
$ws.Range("F4:AO4").Value = $row4

# --- Row 5 ---
$ws.Range("A5").Value = "Colombian Primera B"
$ws.Range("C5").Value = "20:10:00"
$ws.Range("D5").Value = "Real Cartagena"
$ws.Range("E5").Value = "Real Soacha Cundinamarca FC"
$row5 = New-Object 'object[,]' 1,36
$row5[0,0] = 1.65
$row5[0,1] = 1.76
$row5[0,2] = 5.9
$row5[0,3] = 7.8
$row5[0,4] = 3.6
$row5[0,5] = 4
$row5[0,6] = 1.45
$row5[0,7] = 1.09
$row5[0,8] = 3.35
$row5[0,9] = 1.37
$row5[0,10] = 1.76
$row5[0,11] = 2.12
$row5[0,12] = 1.29
$row5[0,13] = 3.9
$row5[0,14] = 2.04
$row5[0,15] = 1.76
$row5[0,16] = 1.15
$row5[0,17] = 2.34
$row5[0,18] = 1000
$row5[0,19] = 1000
$row5[0,20] = 1000
$row5[0,21] = 1000
$row5[0,22] = 29
$row5[0,23] = 42
$row5[0,24] = 1000
$row5[0,25] = 1000
$row5[0,26] = 24
$row5[0,27] = 40
$row5[0,28] = 1000
$row5[0,29] = 1000
$row5[0,30] = 130
$row5[0,31] = 1000
$row5[0,32] = 1000
$row5[0,33] = 1000
$row5[0,34] = 55
$row5[0,35] = 1000
$ws.Range("F5:AO5").Value = $row5

# --- Row 6 ---
$ws.Range("A6").Value = "Brazilian Serie A"
$ws.Range("C6").Value = "20:30:00"
$ws.Range("D6").Value = "Atletico MG"
$ws.Range("E6").Value = "Fortaleza EC"
$row6 = New-Object 'object[,]' 1,36
$row6[0,0] = 1.74
$row6[0,1] = 1.75
$row6[0,2] = 5.9
$row6[0,3] = 6
$row6[0,4] = 3.9
$row6[0,5] = 3.95
$row6[0,6] = 1.44
$row6[0,7] = 1.08
$row6[0,8] = 3.65
$row6[0,9] = 1.35
$row6[0,10] = 1.89
$row6[0,11] = 2.08
$row6[0,12] = 1.33
$row6[0,13] = 3.8
$row6[0,14] = 2.04
$row6[0,15] = 1.92
$row6[0,16] = 1.2
$row6[0,17] = 2.34
$row6[0,18] = 14
$row6[0,19] = 18
$row6[0,20] = 44
$row6[0,21] = 160
$row6[0,22] = 7.8
$row6[0,23] = 8.6
$row6[0,24] = 22
$row6[0,25] = 90
$row6[0,26] = 9
$row6[0,27] = 10.5
$row6[0,28] = 22
$row6[0,29] = 95
$row6[0,30] = 16.5
$row6[0,31] = 19
$row6[0,32] = 42
$row6[0,33] = 140
$row6[0,34] = 12
$row6[0,35] = 110
$ws.Range("F6:AO6").Value = $row6

# --- Row 7 ---
$ws.Range("A7").Value = "Colombian Primera A"
$ws.Range("C7").Value = "22:20:00"
$ws.Range("D7").Value = "Boyaca Chico"
$ws.Range("E7").Value = "Millonarios"
$row7 = New-Object 'object[,]' 1,36
$row7[0,0] = 4.8
$row7[0,1] = 5.4
$row7[0,2] = 1.86
$row7[0,3] = 1.9
$row7[0,4] = 3.5
$row7[0,5] = 3.75
$row7[0,6] = 1.43
$row7[0,7] = 1.08
$row7[0,8] = 3.5
$row7[0,9] = 1.36
$row7[0,10] = 1.84
$row7[0,11] = 2.14
$row7[0,12] = 1.31
$row7[0,13] = 3.85
$row7[0,14] = 1.87
$row7[0,15] = 1.94
$row7[0,16] = 2.1
$row7[0,17] = 1.23
$row7[0,18] = 980
$row7[0,19] = 9.2
$row7[0,20] = 980
$row7[0,21] = 980
$row7[0,22] = 980
$row7[0,23] = 9.4
$row7[0,24] = 980
$row7[0,25] = 980
$row7[0,26] = 980
$row7[0,27] = 980
$row7[0,28] = 980
$row7[0,29] = 55
$row7[0,30] = 140
$row7[0,31] = 80
$row7[0,32] = 100
$row7[0,33] = 1000
$row7[0,34] = 1000
$row7[0,35] = 980
$ws.Range("F7:AO7").Value = $row7

